# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 246.9852506941017
$ws.Range("G2").Value = 265.7776617428604

$ws.Range("B3").Value = [double]"1.332267629550188e-15"
$ws.Range("C3").Value = [double]"2.220651329265522e-06"
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 7.194163148936941

$ws.Range("B4").Value = 3.182878228561681
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 28.30127388105354

$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 86.29678392075563
$ws.Range("D5").Value = 0.7127328510149897
$ws.Range("E5").Value = 246.9852506941017
$ws.Range("G5").Value = 334.7234868868072

$ws.Range("B6").Value = 0.02258322285507441
$ws.Range("C6").Value = 0.3375848360084654
$ws.Range("D6").Value = 3.082599426703578
$ws.Range("E6").Value = 246.9852506941017
$ws.Range("G6").Value = 250.4280181796688
